$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.569.48"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.920.98"
$ws.Range("E3").Value = "  +3.27%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'247.32"
$ws.Range("E5").Value = "  +4.70%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").Value = "'0.2880"
$ws.Range("E8").Value = "  +3.19%  "
$ws.Range("D9").Value = "'0.06808"
$ws.Range("E9").Value = "  +6.06%  "
$ws.Range("D10").Value = "'105.36"
$ws.Range("E10").Value = "  +5.33%  "
$ws.Range("D11").Value = "'18.37"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "1.917.46"
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("D13").Value = "'0.07688"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "'5.293"
$ws.Range("E14").Value = "  +5.81%  "
$ws.Range("D15").Value = "'0.6678"
$ws.Range("E15").Value = "  +6.00%  "
$ws.Range("D16").Value = "'288.92"
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("D17").Value = "30.583.74"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007592"
$ws.Range("E19").Value = "  +2.62%  "
$ws.Range("D20").Value = "'12.92"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("D21").Value = "'5.516"
$ws.Range("E21").Value = "  +9.83%  "
$ws.Range("D22").Value = "2.169.40"
$ws.Range("E22").Value = "  +3.66%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'6.302"
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("D25").Value = "'9.378"
$ws.Range("E25").Value = "  +2.97%  "
$ws.Range("D26").Value = "'168.75"
$ws.Range("E26").Value = "  +2.27%  "
$ws.Range("D27").Value = "'21.36"
$ws.Range("E27").Value = "  +10.24%  "
$ws.Range("D28").Value = "'2.122"
$ws.Range("E28").Value = "  +9.58%  "
$ws.Range("D29").Value = "'0.1071"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").Value = "'1.396"
$ws.Range("E30").Value = "  +5.59%  "
$ws.Range("D31").Value = "'4.172"
$ws.Range("E31").Value = "  +4.14%  "
$ws.Range("D32").Value = "'4.085"
$ws.Range("E32").Value = "  +6.04%  "
$ws.Range("D33").Value = "'0.05045"
$ws.Range("E33").Value = "  +3.34%  "
$ws.Range("D34").Value = "'0.7372"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'1.148"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").Value = "'0.02074"
$ws.Range("E36").Value = "  +8.91%  "
$ws.Range("D37").Value = "'2.751"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").Value = "'2.690"
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").Value = "'2.065"
$ws.Range("E39").Value = "  +4.74%  "
$ws.Range("D40").Value = "'111.44"
$ws.Range("E40").Value = "  +4.14%  "
$ws.Range("D41").Value = "'0.8796"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("D42").Value = "'0.4406"
$ws.Range("E42").Value = "  +8.48%  "
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "'67.23"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("D46").Value = "'7.271"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("D47").Value = "'9.339"
$ws.Range("E47").Value = "  +2.97%  "
$ws.Range("D48").Value = "'48.66"
$ws.Range("E48").Value = "  +18.42%  "
$ws.Range("D49").Value = "'0.1232"
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.480"
$ws.Range("E50").Value = "  +10.93%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.4100"
$ws.Range("E51").Value = "  +9.97%  "
